$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new "Wins" / "Losses" / "Ties" columns, matching the existing
# header style (bold, centered, bordered) by copying format from A1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-55: every player on this roster shares the same team record.
$lastRow = 55
$wins = 74
$losses = 88
$ties = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
